$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "[1,2]"
$ws.Range("D3").Value = "['C3','C4']"
$ws.Range("F2").Value = "[1.0, 0.4]"
$ws.Range("F3").Value = "[2.0]"
$ws.Range("G2").Value = "[0.1, 0.2]"
$ws.Range("G3").Value = "[0.2]"

$ws.Range("H14").Select()
